# Sudip - Time Sheet (Week 9): record hours 8 & 9 work entries and move
# the active selection to B9, matching the author's commit
# "Edited some rows of time sheet."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$note = "Worked on redirecting non-logged in user to login page when click on Appointment tab "

$ws.Range("B9").Value = $note
$ws.Range("B10").Value = $note

$ws.Range("B9").Select()
